$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value would be parsed as a plain number by Excel
# (losing formatting such as trailing zeros) need to be forced to Text format
# before the value is assigned, so they remain stored as strings, matching the source data.
$textFormatCells = @("D5", "D6", "D9", "D10", "D12", "D15", "D19", "D20", "D21", "D22", "D24", "D25", "D26", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D42", "D43", "D44", "D47", "D49", "D50", "D51")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values scraped from the source for this run.
$ws.Range('D2').Value = '55.231.34'
$ws.Range('D3').Value = '2.939.90'
$ws.Range('E3').Value = '  -6.91%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '483.00'
$ws.Range('E5').Value = '  -8.11%  '
$ws.Range('D6').Value = '129.54'
$ws.Range('E6').Value = '  -2.76%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '2.939.90'
$ws.Range('E8').Value = '  -6.94%  '
$ws.Range('D9').Value = '0.415'
$ws.Range('E9').Value = '  -8.12%  '
$ws.Range('D10').Value = '6.94'
$ws.Range('E10').Value = '  -4.95%  '
$ws.Range('E11').Value = '  -10.45%  '
$ws.Range('D12').Value = '0.345'
$ws.Range('E12').Value = '  -11.67%  '
$ws.Range('E13').Value = '  -1.05%  '
$ws.Range('D14').Value = '3.443.12'
$ws.Range('E14').Value = '  -7.01%  '
$ws.Range('D15').Value = '24.11'
$ws.Range('E15').Value = '  -6.77%  '
$ws.Range('D16').Value = '55.172.19'
$ws.Range('E16').Value = '  -4.53%  '
$ws.Range('D17').Value = '2.932.75'
$ws.Range('E17').Value = '  -7.24%  '
$ws.Range('E18').Value = '  -9.60%  '
$ws.Range('D19').Value = '5.54'
$ws.Range('E19').Value = '  -4.82%  '
$ws.Range('D20').Value = '11.89'
$ws.Range('E20').Value = '  -9.00%  '
$ws.Range('D21').Value = '7.36'
$ws.Range('E21').Value = '  -8.55%  '
$ws.Range('D22').Value = '308.79'
$ws.Range('E22').Value = '  -10.87%  '
$ws.Range('E23').Value = '  +0.16%  '
$ws.Range('D24').Value = '0.457'
$ws.Range('E24').Value = '  -10.75%  '
$ws.Range('D25').Value = '59.49'
$ws.Range('E25').Value = '  -14.58%  '
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.16%  '
$ws.Range('E27').Value = '  -6.25%  '
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('D29').Value = '0.0₃0839'
$ws.Range('E29').Value = '  -12.61%  '
$ws.Range('D30').Value = '6.49'
$ws.Range('E30').Value = '  -5.18%  '
$ws.Range('D31').Value = '1.17'
$ws.Range('E31').Value = '  -3.37%  '
$ws.Range('D32').Value = '6.41'
$ws.Range('E32').Value = '  -7.25%  '
$ws.Range('D33').Value = '1.67'
$ws.Range('E33').Value = '  -11.46%  '
$ws.Range('D34').Value = '19.13'
$ws.Range('E34').Value = '  -11.92%  '
$ws.Range('D35').Value = '146.39'
$ws.Range('E35').Value = '  -8.58%  '
$ws.Range('D36').Value = '4.31'
$ws.Range('E36').Value = '  -11.95%  '
$ws.Range('D37').Value = '5.59'
$ws.Range('E37').Value = '  -10.48%  '
$ws.Range('E38').Value = '  -10.58%  '
$ws.Range('D39').Value = '23.65'
$ws.Range('E39').Value = '  -8.31%  '
$ws.Range('D40').Value = '0.0640'
$ws.Range('E40').Value = '  -8.14%  '
$ws.Range('D41').Value = '2.969.35'
$ws.Range('E41').Value = '  -6.80%  '
$ws.Range('D42').Value = '0.998'
$ws.Range('E42').Value = '  -0.19%  '
$ws.Range('D43').Value = '35.76'
$ws.Range('E43').Value = '  -11.96%  '
$ws.Range('D44').Value = '0.987'
$ws.Range('E44').Value = '  -8.76%  '
$ws.Range('E45').Value = '  -10.92%  '
$ws.Range('E46').Value = '  -7.53%  '
$ws.Range('D47').Value = '3.49'
$ws.Range('E47').Value = '  -11.51%  '
$ws.Range('D48').Value = '2.107.40'
$ws.Range('E48').Value = '  -7.13%  '
$ws.Range('D49').Value = '0.0227'
$ws.Range('E49').Value = '  -3.66%  '
$ws.Range('D50').Value = '18.71'
$ws.Range('E50').Value = '  -8.83%  '
$ws.Range('D51').Value = '5.50'
$ws.Range('E51').Value = '  -11.04%  '
